# R&D_2025_Combined_Frequency_Tables.xlsx — add Scores_Correlation sheet,
# refresh header styling across sheets, re-activate the first tab, and
# normalize All_Frequencies' page margins to Excel's modern defaults.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add the new "Scores_Correlation" worksheet as the last tab
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws7 = $wb.Worksheets.Add($null, $lastSheet)
$ws7.Name = "Scores_Correlation"

$ws7.Range("A1").Value = "Variable"
$ws7.Range("B1").Value = "Variable"
$ws7.Range("C1").Value = "knowledge_score"
$ws7.Range("D1").Value = "awareness_Score"
$ws7.Range("E1").Value = "perception_Score"

$ws7.Range("A2").Value = 0
$ws7.Range("B2").Value = "knowledge_score"
$ws7.Range("C2").Value = "1.000 (p=0.000)"
$ws7.Range("D2").Value = "0.321 (p=0.000)"
$ws7.Range("E2").Value = "0.219 (p=0.001)"

$ws7.Range("A3").Value = 1
$ws7.Range("B3").Value = "awareness_Score"
$ws7.Range("C3").Value = "0.321 (p=0.000)"
$ws7.Range("D3").Value = "1.000 (p=0.000)"
$ws7.Range("E3").Value = "0.239 (p=0.000)"

$ws7.Range("A4").Value = 2
$ws7.Range("B4").Value = "perception_Score"
$ws7.Range("C4").Value = "0.219 (p=0.001)"
$ws7.Range("D4").Value = "0.239 (p=0.000)"
$ws7.Range("E4").Value = "1.000 (p=0.000)"

# Header formatting: bold font, thin box border, centered/top aligned —
# matches the look already used by every other table's header row.
$hdr7 = $ws7.Range("A1:E1")
$hdr7.Font.Bold = $true
$hdr7.Borders.LineStyle = 1
$hdr7.HorizontalAlignment = -4108
$hdr7.VerticalAlignment = -4160

# ---------------------------------------------------------------------
# 2) Refresh the header-row styling on the other data sheets so every
#    header cell shares the same (freshly rebuilt) style definition.
# ---------------------------------------------------------------------
$styleSrc = $ws7.Range("A1")
$styleSrc.Copy()

$wb.Worksheets.Item("Frequency_Tables").Range("A1:D1").PasteSpecial(-4122)
$wb.Worksheets.Item("Awareness_Perception_Frequencies").Range("A1:F1").PasteSpecial(-4122)
$wb.Worksheets.Item("Scores_Statistics").Range("B1:C1").PasteSpecial(-4122)
$wb.Worksheets.Item("Scores_Statistics").Range("A2:A4").PasteSpecial(-4122)
$wb.Worksheets.Item("Score_Categories").Range("B1:D1").PasteSpecial(-4122)
$wb.Worksheets.Item("Score_Categories").Range("A2:A3").PasteSpecial(-4122)
$wb.Worksheets.Item("Score_Percentages").Range("B1:D1").PasteSpecial(-4122)
$wb.Worksheets.Item("Score_Percentages").Range("A2:A3").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) All_Frequencies: normalize page margins to Excel's modern defaults
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("All_Frequencies")
$ws1.PageSetup.LeftMargin = 50.4
$ws1.PageSetup.RightMargin = 50.4
$ws1.PageSetup.TopMargin = 54
$ws1.PageSetup.BottomMargin = 54
$ws1.PageSetup.HeaderMargin = 21.599999999999998
$ws1.PageSetup.FooterMargin = 21.599999999999998

# ---------------------------------------------------------------------
# 4) Re-select the first tab so it (not the newly-added sheet) is the
#    one marked active/selected when the workbook is reopened.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1").Select() | Out-Null
